# Updated cryptos list on Fri May 10 02:58:16 UTC 2024 with GitHub Actions
#
# Note: every assigned literal below begins with an escaped leading apostrophe
# ('' at the start of each single-quoted PowerShell string, i.e. the value
# actually assigned starts with a single "'" character). That forces Excel to
# store the value as plain text, exactly like the source workbook's inline
# strings, instead of auto-converting number-looking text (e.g. "592.18",
# "132.00", "0.0000234") into numeric cell values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''62.841.95'
$ws.Range("E2").Value = '''  +1.92%  '

$ws.Range("D3").Value = '''3.032.50'
$ws.Range("E3").Value = '''  +1.10%  '

$ws.Range("E4").Value = '''  +0.10%  '

$ws.Range("D5").Value = '''592.18'
$ws.Range("E5").Value = '''  -1.07%  '

$ws.Range("D6").Value = '''153.04'
$ws.Range("E6").Value = '''  +6.02%  '

$ws.Range("E7").Value = '''  +0.05%  '

$ws.Range("D8").Value = '''3.026.38'
$ws.Range("E8").Value = '''  +0.91%  '

$ws.Range("D9").Value = '''0.518'
$ws.Range("E9").Value = '''  -0.70%  '

$ws.Range("D10").Value = '''6.38'
$ws.Range("E10").Value = '''  +8.17%  '

$ws.Range("D11").Value = '''0.150'
$ws.Range("E11").Value = '''  +2.12%  '

$ws.Range("D12").Value = '''0.464'
$ws.Range("E12").Value = '''  +0.29%  '

$ws.Range("D13").Value = '''0.0000234'
$ws.Range("E13").Value = '''  +2.50%  '

$ws.Range("D14").Value = '''35.51'
$ws.Range("E14").Value = '''  +3.39%  '

$ws.Range("E15").Value = '''  +1.88%  '

$ws.Range("D16").Value = '''3.536.46'
$ws.Range("E16").Value = '''  +1.19%  '

$ws.Range("E17").Value = '''  +0.99%  '

$ws.Range("D18").Value = '''62.864.49'
$ws.Range("E18").Value = '''  +2.06%  '

$ws.Range("D19").Value = '''3.032.39'
$ws.Range("E19").Value = '''  +1.34%  '

$ws.Range("D20").Value = '''451.44'
$ws.Range("E20").Value = '''  -0.82%  '

$ws.Range("D21").Value = '''14.29'
$ws.Range("E21").Value = '''  +1.86%  '

$ws.Range("E22").Value = '''  +0.77%  '

$ws.Range("D23").Value = '''7.47'
$ws.Range("E23").Value = '''  +1.42%  '

$ws.Range("D24").Value = '''83.12'
$ws.Range("E24").Value = '''  +0.86%  '

$ws.Range("D25").Value = '''2.31'
$ws.Range("E25").Value = '''  +4.60%  '

$ws.Range("D26").Value = '''11.06'
$ws.Range("E26").Value = '''  +5.67%  '

$ws.Range("D27").Value = '''12.23'
$ws.Range("E27").Value = '''  +0.50%  '

$ws.Range("E28").Value = '''  -0.06%  '

$ws.Range("D29").Value = '''7.56'
$ws.Range("E29").Value = '''  +7.79%  '

$ws.Range("D30").Value = '''2.69'
$ws.Range("E30").Value = '''  +0.63%  '

$ws.Range("E31").Value = '''  +7.72%  '

$ws.Range("E32").Value = '''  +0.13%  '

$ws.Range("D33").Value = '''27.52'
$ws.Range("E33").Value = '''  +0.52%  '

$ws.Range("E34").Value = '''  +2.82%  '

$ws.Range("D35").Value = '''0.0₃0867'
$ws.Range("E35").Value = '''  +5.91%  '

$ws.Range("E36").Value = '''  +1.61%  '

$ws.Range("D37").Value = '''5.91'
$ws.Range("E37").Value = '''  +2.67%  '

$ws.Range("D38").Value = '''3.15'
$ws.Range("E38").Value = '''  +8.81%  '

$ws.Range("D39").Value = '''2.09'
$ws.Range("E39").Value = '''  +0.24%  '

$ws.Range("D40").Value = '''50.53'
$ws.Range("E40").Value = '''  +0.44%  '

$ws.Range("B41").Value = '''Kaspa'
$ws.Range("C41").Value = '''https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").Value = '''0.127'
$ws.Range("E41").Value = '''  +3.30%  '

$ws.Range("B42").Value = '''Cosmos'
$ws.Range("C42").Value = '''https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D42").Value = '''9.09'
$ws.Range("E42").Value = '''  -1.12%  '

$ws.Range("E43").Value = '''  +16.36%  '

$ws.Range("D44").Value = '''42.47'
$ws.Range("E44").Value = '''  +8.07%  '

$ws.Range("D45").Value = '''395.60'
$ws.Range("E45").Value = '''  -1.19%  '

$ws.Range("D46").Value = '''0.0359'
$ws.Range("E46").Value = '''  +1.56%  '

$ws.Range("D47").Value = '''2.737.69'
$ws.Range("E47").Value = '''  +0.65%  '

$ws.Range("D48").Value = '''132.00'
$ws.Range("E48").Value = '''  -0.95%  '

$ws.Range("E49").Value = '''  +0.03%  '

$ws.Range("E50").Value = '''  +3.22%  '

$ws.Range("D51").Value = '''24.23'
$ws.Range("E51").Value = '''  +3.11%  '
